{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2022-11-21 Monday\", \"2022-11-22 Tuesday\"],\n  [\"8+36=\", \"56-0=\"],\n  [\"46+15=\", \"76+20=\"],\n  [\"89-42=\", \"55-35=\"],\n  [\"64+1=\", \"63-35=\"],\n  [\"74-39=\", \"18-8=\"],\n  [\"17-11=\", \"91-75=\"],\n  [\"74+2=\", \"41+34=\"],\n  [\"78-46=\", \"33+51=\"],\n  [\"19+61=\", \"11+54=\"],\n  [\"70-64=\", \"5+19=\"],\n  [\"47+27=\", \"84-55=\"],\n  [\"52+9=\", \"19+13=\"],\n  [\"97-20=\", \"89-62=\"],\n  [\"67-54=\", \"29+10=\"],\n  [\"37+51=\", \"31-13=\"],\n  [\"79-52=\", \"49-7=\"],\n  [\"94-12=\", \"1+80=\"],\n  [\"1+8=\", \"24-2=\"],\n  [\"58-42=\", \"98-90=\"],\n  [\"82-71=\", \"93-51=\"],\n  [\"95-87=\", \"15+43=\"],\n  [\"6+42=\", \"10+11=\"],\n  [\"7+49=\", \"2+92=\"],\n  [\"58+12=\", \"10+10=\"],\n  [\"1+50=\", \"6+87=\"],\n  [\"20-4=\", \"83+8=\"],\n  [\"58+29=\", \"59+0=\"],\n  [\"33+55=\", \"45+46=\"],\n  [\"44-20=\", \"61+3=\"],\n  [\"9+3=\", \"91+4=\"],\n  [\"15-5=\", \"10+74=\"],\n  [\"96-72=\", \"15+43=\"],\n  [\"31+12=\", \"71-4=\"],\n  [\"41+24=\", \"43-3=\"],\n  [\"22+5=\", \"29+43=\"],\n  [\"82+3=\", \"36-15=\"],\n  [\"4+18=\", \"87-66=\"],\n  [\"13+65=\", \"40+7=\"],\n  [\"9+83=\", \"21+11=\"],\n  [\"52-3=\", \"93-68=\"],\n  [\"17+76=\", \"65-60=\"],\n  [\"48+24=\", \"26+44=\"],\n  [\"62+20=\", \"92-64=\"],\n  [\"61+11=\", \"3+39=\"],\n  [\"22-6=\", \"93-57=\"],\n  [\"87-40=\", \"92-60=\"],\n  [\"9+17=\", \"87+9=\"],\n  [\"24+47=\", \"98-48=\"],\n  [\"43-35=\", \"40+43=\"],\n  [\"40+24=\", \"26+62=\"],\n  [\"26-7=\", \"45+26=\"],\n  [\"62-47=\", \"91-36=\"],\n  [\"39+12=\", \"95-83=\"],\n  [\"58-45=\", \"9-7=\"],\n  [\"75-8=\", \"24-22=\"],\n  [\"12+45=\", \"57-6=\"],\n  [\"3+47=\", \"34+52=\"],\n  [\"93-17=\", \"35+63=\"],\n  [\"94-75=\", \"3+2=\"],\n  [\"28+4=\", \"5+71=\"],\n  [\"17+62=\", \"61+6=\"],\n  [\"27+17=\", \"61-18=\"],\n  [\"12+75=\", \"80+7=\"],\n  [\"97-38=\", \"35+0=\"],\n  [\"0+60=\", \"29+52=\"],\n  [\"57-41=\", \"40-20=\"],\n  [\"65-1=\", \"95-69=\"],\n  [\"97-12=\", \"72+5=\"],\n  [\"0+66=\", \"17+33=\"],\n  [\"21+55=\", \"50+4=\"],\n  [\"21+45=\", \"79-35=\"],\n  [\"17+69=\", \"26-12=\"],\n  [\"93-56=\", \"88-32=\"],\n  [\"22+6=\", \"22+14=\"],\n  [\"67-26=\", \"8+62=\"],\n  [\"17+3=\", \"63-22=\"],\n  [\"27-8=\", \"64+29=\"],\n  [\"32-25=\", \"10+7=\"],\n  [\"78+9=\", \"60-46=\"],\n  [\"74-35=\", \"61+9=\"],\n  [\"60+3=\", \"74-55=\"],\n  [\"26+71=\", \"41+42=\"],\n  [\"30+50=\", \"11+4=\"],\n  [\"18+72=\", \"85-25=\"],\n  [\"14+37=\", \"74-15=\"],\n  [\"1+64=\", \"48+33=\"],\n  [\"51+45=\", \"96-37=\"],\n  [\"20-6=\", \"95-32=\"],\n  [\"24+34=\", \"17+0=\"],\n  [\"49-27=\", \"91-24=\"],\n  [\"86+2=\", \"68-50=\"],\n  [\"20+77=\", \"52+42=\"],\n  [\"8+39=\", \"76+1=\"],\n  [\"88-15=\", \"6+93=\"],\n  [\"93-84=\", \"50-6=\"],\n  [\"8-6=\", \"76-18=\"],\n  [\"86+3=\", \"56+16=\"],\n  [\"34-13=\", \"57+4=\"],\n  [\"62-25=\", \"76+18=\"],\n  [\"11+40=\", \"94-47=\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + findText);\n  }\n  for (const item of results.items) {\n    item.insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n    @(\"2022-11-21 Monday\", \"2022-11-22 Tuesday\"),\n    @(\"8+36=\", \"56-0=\"),\n    @(\"46+15=\", \"76+20=\"),\n    @(\"89-42=\", \"55-35=\"),\n    @(\"64+1=\", \"63-35=\"),\n    @(\"74-39=\", \"18-8=\"),\n    @(\"17-11=\", \"91-75=\"),\n    @(\"74+2=\", \"41+34=\"),\n    @(\"78-46=\", \"33+51=\"),\n    @(\"19+61=\", \"11+54=\"),\n    @(\"70-64=\", \"5+19=\"),\n    @(\"47+27=\", \"84-55=\"),\n    @(\"52+9=\", \"19+13=\"),\n    @(\"97-20=\", \"89-62=\"),\n    @(\"67-54=\", \"29+10=\"),\n    @(\"37+51=\", \"31-13=\"),\n    @(\"79-52=\", \"49-7=\"),\n    @(\"94-12=\", \"1+80=\"),\n    @(\"1+8=\", \"24-2=\"),\n    @(\"58-42=\", \"98-90=\"),\n    @(\"82-71=\", \"93-51=\"),\n    @(\"95-87=\", \"15+43=\"),\n    @(\"6+42=\", \"10+11=\"),\n    @(\"7+49=\", \"2+92=\"),\n    @(\"58+12=\", \"10+10=\"),\n    @(\"1+50=\", \"6+87=\"),\n    @(\"20-4=\", \"83+8=\"),\n    @(\"58+29=\", \"59+0=\"),\n    @(\"33+55=\", \"45+46=\"),\n    @(\"44-20=\", \"61+3=\"),\n    @(\"9+3=\", \"91+4=\"),\n    @(\"15-5=\", \"10+74=\"),\n    @(\"96-72=\", \"15+43=\"),\n    @(\"31+12=\", \"71-4=\"),\n    @(\"41+24=\", \"43-3=\"),\n    @(\"22+5=\", \"29+43=\"),\n    @(\"82+3=\", \"36-15=\"),\n    @(\"4+18=\", \"87-66=\"),\n    @(\"13+65=\", \"40+7=\"),\n    @(\"9+83=\", \"21+11=\"),\n    @(\"52-3=\", \"93-68=\"),\n    @(\"17+76=\", \"65-60=\"),\n    @(\"48+24=\", \"26+44=\"),\n    @(\"62+20=\", \"92-64=\"),\n    @(\"61+11=\", \"3+39=\"),\n    @(\"22-6=\", \"93-57=\"),\n    @(\"87-40=\", \"92-60=\"),\n    @(\"9+17=\", \"87+9=\"),\n    @(\"24+47=\", \"98-48=\"),\n    @(\"43-35=\", \"40+43=\"),\n    @(\"40+24=\", \"26+62=\"),\n    @(\"26-7=\", \"45+26=\"),\n    @(\"62-47=\", \"91-36=\"),\n    @(\"39+12=\", \"95-83=\"),\n    @(\"58-45=\", \"9-7=\"),\n    @(\"75-8=\", \"24-22=\"),\n    @(\"12+45=\", \"57-6=\"),\n    @(\"3+47=\", \"34+52=\"),\n    @(\"93-17=\", \"35+63=\"),\n    @(\"94-75=\", \"3+2=\"),\n    @(\"28+4=\", \"5+71=\"),\n    @(\"17+62=\", \"61+6=\"),\n    @(\"27+17=\", \"61-18=\"),\n    @(\"12+75=\", \"80+7=\"),\n    @(\"97-38=\", \"35+0=\"),\n    @(\"0+60=\", \"29+52=\"),\n    @(\"57-41=\", \"40-20=\"),\n    @(\"65-1=\", \"95-69=\"),\n    @(\"97-12=\", \"72+5=\"),\n    @(\"0+66=\", \"17+33=\"),\n    @(\"21+55=\", \"50+4=\"),\n    @(\"21+45=\", \"79-35=\"),\n    @(\"17+69=\", \"26-12=\"),\n    @(\"93-56=\", \"88-32=\"),\n    @(\"22+6=\", \"22+14=\"),\n    @(\"67-26=\", \"8+62=\"),\n    @(\"17+3=\", \"63-22=\"),\n    @(\"27-8=\", \"64+29=\"),\n    @(\"32-25=\", \"10+7=\"),\n    @(\"78+9=\", \"60-46=\"),\n    @(\"74-35=\", \"61+9=\"),\n    @(\"60+3=\", \"74-55=\"),\n    @(\"26+71=\", \"41+42=\"),\n    @(\"30+50=\", \"11+4=\"),\n    @(\"18+72=\", \"85-25=\"),\n    @(\"14+37=\", \"74-15=\"),\n    @(\"1+64=\", \"48+33=\"),\n    @(\"51+45=\", \"96-37=\"),\n    @(\"20-6=\", \"95-32=\"),\n    @(\"24+34=\", \"17+0=\"),\n    @(\"49-27=\", \"91-24=\"),\n    @(\"86+2=\", \"68-50=\"),\n    @(\"20+77=\", \"52+42=\"),\n    @(\"8+39=\", \"76+1=\"),\n    @(\"88-15=\", \"6+93=\"),\n    @(\"93-84=\", \"50-6=\"),\n    @(\"8-6=\", \"76-18=\"),\n    @(\"86+3=\", \"56+16=\"),\n    @(\"34-13=\", \"57+4=\"),\n    @(\"62-25=\", \"76+18=\"),\n    @(\"11+40=\", \"94-47=\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $result = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $wdReplaceAll)\n    if (-not $result) {\n        Write-Output \"WARNING: replace failed for $findText\"\n    }\n}\n"}
